$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 6823
$ws.Cells.Item(3, 6).Value = 836
$ws.Cells.Item(8, 6).Value = 25
$ws.Cells.Item(11, 6).Value = 1125
$ws.Cells.Item(12, 6).Value = 883
$ws.Cells.Item(13, 6).Value = 17
$ws.Cells.Item(14, 6).Value = 710
$ws.Cells.Item(15, 6).Value = 1023
$ws.Cells.Item(16, 6).Value = 1383
$ws.Cells.Item(17, 6).Value = 52
$ws.Cells.Item(19, 6).Value = 1556
$ws.Cells.Item(21, 6).Value = 586
$ws.Cells.Item(23, 6).Value = 6
$ws.Cells.Item(25, 6).Value = 1073
$ws.Cells.Item(26, 6).Value = 1504
$ws.Cells.Item(27, 6).Value = 735
$ws.Cells.Item(28, 6).Value = 577
$ws.Cells.Item(30, 6).Value = 463
$ws.Cells.Item(32, 6).Value = 1014
$ws.Cells.Item(34, 6).Value = 288
$ws.Cells.Item(35, 6).Value = 2384
$ws.Cells.Item(36, 6).Value = 267
$ws.Cells.Item(37, 6).Value = 1307
$ws.Cells.Item(40, 6).Value = 3908
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(25, 6).Value = 232
$ws.Cells.Item(28, 6).Value = 54
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 94
$ws.Cells.Item(5, 6).Value = 1652
$ws.Cells.Item(6, 6).Value = 449
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 94
$ws.Cells.Item(4, 6).Value = 1652
$ws.Cells.Item(5, 6).Value = 449
$ws.Cells.Item(8, 6).Value = 6823
$ws.Cells.Item(9, 6).Value = 836
$ws.Cells.Item(14, 6).Value = 25
$ws.Cells.Item(17, 6).Value = 1125
$ws.Cells.Item(18, 6).Value = 883
$ws.Cells.Item(19, 6).Value = 710
$ws.Cells.Item(23, 6).Value = 1023
$ws.Cells.Item(24, 6).Value = 1383
$ws.Cells.Item(25, 6).Value = 52
$ws.Cells.Item(27, 6).Value = 1556
$ws.Cells.Item(29, 6).Value = 586
$ws.Cells.Item(32, 6).Value = 1073
$ws.Cells.Item(33, 6).Value = 1504
$ws.Cells.Item(34, 6).Value = 735
$ws.Cells.Item(35, 6).Value = 577
$ws.Cells.Item(37, 6).Value = 463
$ws.Cells.Item(41, 6).Value = 1014
$ws.Cells.Item(43, 6).Value = 288
$ws.Cells.Item(44, 6).Value = 2384
$ws.Cells.Item(45, 6).Value = 232
$ws.Cells.Item(49, 6).Value = 1307
$ws.Cells.Item(51, 6).Value = 3908
